$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value2 = "[Emily%Mosites%NULL%1,             Erin M.%Parker%NULL%2,             Kristie E. N.%Clarke%NULL%1,             Jessie M.%Gaeta%NULL%1,             Travis P.%Baggett%NULL%1,             Elizabeth%Imbert%NULL%1,             Madeline%Sankaran%NULL%2,             Ashley%Scarborough%NULL%2,             Karin%Huster%NULL%2,             Matt%Hanson%NULL%1,             Elysia%Gonzales%NULL%2,             Jody%Rauch%NULL%2,             Libby%Page%NULL%1,             Temet M.%McMichael%NULL%2,             Ryan%Keating%NULL%2,             Grace E.%Marx%NULL%2,             Tom%Andrews%NULL%2,             Kristine%Schmit%NULL%2,             Sapna Bamrah%Morris%NULL%2,             Nicole F.%Dowling%NULL%2,             Georgina%Peacock%NULL%2,             NULL%NULL%NULL%5,             NULL%NULL%NULL%0,             NULL%NULL%NULL%0,             Ann%Buff%NULL%1,             Calla%Jamison%NULL%2,             Ruthanne%Marcus%NULL%1,             Carol Y.%Rao%NULL%2,             Julie L.%Self%NULL%2,             Farrell%Tobolowsky%NULL%1,             Samantha%Williams%NULL%1,             Meagan%Kay%NULL%2,             Naveena%Bobba%NULL%1,             Stephanie%Cohen%NULL%2,             Jonathan%Fuchs%NULL%1,             Trang%Nguyen%NULL%1,             Julie%Stoltey%NULL%1]"
$ws.Cells.Item(3, 5).Value2 = "[Travis P.%Baggett%xref no email%2,        Harrison%Keyes%xref no email%1,        Nora%Sporn%xref no email%1,        Jessie M.%Gaeta%xref no email%2]"
$ws.Cells.Item(4, 5).Value2 = "[Isaac%Ghinai%isaac.ghinai@cityofchicago.org%1,             Elizabeth S%Davis%NULL%2,             Elizabeth S%Davis%NULL%0,             Stockton%Mayer%NULL%1,             Karrie-Ann%Toews%NULL%1,             Thomas D%Huggett%NULL%1,             Nyssa%Snow-Hill%NULL%1,             Omar%Perez%NULL%1,             Mary K%Hayden%NULL%1,             Seena%Tehrani%NULL%1,             A Justine%Landi%NULL%1,             Stephanie%Crane%NULL%1,             Elizabeth%Bell%NULL%1,             Joy-Marie%Hermes%NULL%1,             Kush%Desai%NULL%1,             Michelle%Godbee%NULL%1,             Naman%Jhaveri%NULL%1,             Brian%Borah%NULL%1,             Tracy%Cable%NULL%1,             Sofia%Sami%NULL%1,             Laura%Nozicka%NULL%1,             Yi-Shin%Chang%NULL%1,             Aditi%Jagadish%NULL%1,             Mark%Chee%NULL%1,             Brynna%Thigpen%NULL%1,             Christopher%Llerena%NULL%1,             Minh%Tran%NULL%1,             Divya Meher%Surabhi%NULL%1,             Emilia D%Smith%NULL%1,             Rosemary G%Remus%NULL%1,             Roweine%Staszcuk%NULL%1,             Evelyn%Figueroa%NULL%1,             Paul%Leo%NULL%1,             Wayne M%Detmer%NULL%1,             Evan%Lyon%NULL%1,             Sarah%Carreon%NULL%1,             Stacey%Hoferka%NULL%1,             Kathleen A%Ritger%NULL%1,             Wilnise%Jasmin%NULL%1,             Prathima%Nagireddy%NULL%1,             Jennifer Y%Seo%NULL%1,             Marielle J%Fricchione%NULL%1,             Janna L%Kerins%NULL%1,             Stephanie R%Black%NULL%1,             Lisa Morrison%Butler%NULL%1,             Kimberly%Howard%NULL%1,             Maura%McCauley%NULL%1,             Todd%Fraley%NULL%1,             M Allison%Arwady%NULL%1,             Stephanie%Gretsch%NULL%1,             Megan%Cunningham%NULL%1,             Massimo%Pacilli%NULL%1,             Peter S%Ruestow%NULL%2,             Peter S%Ruestow%NULL%0,             Emily%Mosites%NULL%3,             Elizabeth%Avery%NULL%1,             Joshua%Longcoy%NULL%1,             Elizabeth B%Lynch%NULL%1,             Jennifer E%Layden%NULL%1]"
$ws.Cells.Item(5, 5).Value2 = "[Matthew D.%Hickey%xref no email%1,        Elizabeth%Imbert%xref no email%1,        David V.%Glidden%xref no email%2,        Jan Bing%Del Rosario%xref no email%1,        Mary%Chong%xref no email%1,        Angelo%Clemenzi-Allen%xref no email%1,        Jon%Oskarsson%xref no email%1,        Elise D.%Riley%xref no email%1,        Monica%Gandhi%xref no email%2,        Diane V.%Havlir%xref no email%1]"
$ws.Cells.Item(6, 5).Value2 = "[Emeline%Han%NULL%1,             Melisa Mei Jin%Tan%NULL%1,             Eva%Turk%NULL%1,             Devi%Sridhar%NULL%1,             Gabriel M%Leung%NULL%1,             Kenji%Shibuya%NULL%1,             Nima%Asgari%NULL%1,             Juhwan%Oh%NULL%1,             Alberto L%García-Basteiro%NULL%1,             Johanna%Hanefeld%NULL%1,             Alex R%Cook%NULL%1,             Li Yang%Hsu%NULL%1,             Yik Ying%Teo%NULL%1,             David%Heymann%NULL%1,             Helen%Clark%NULL%1,             Martin%McKee%NULL%1,             Helena%Legido-Quigley%NULL%1]"
$ws.Cells.Item(7, 5).Value2 = "[Heather E.%Hsu%NULL%1,             Erin M.%Ashe%NULL%1,             Michael%Silverstein%NULL%1,             Melissa%Hofman%NULL%1,             Samantha J.%Lange%NULL%1,             Hilda%Razzaghi%NULL%1,             Rebecca G.%Mishuris%NULL%1,             Ravin%Davidoff%NULL%1,             Erin M.%Parker%NULL%0,             Ana%Penman-Aguilar%NULL%1,             Kristie E.N.%Clarke%NULL%1,             Anna%Goldman%NULL%1,             Thea L.%James%NULL%1,             Karen%Jacobson%NULL%1,             Karen E.%Lasser%NULL%1,             Ziming%Xuan%NULL%1,             Georgina%Peacock%NULL%0,             Nicole F.%Dowling%NULL%0,             Alyson B.%Goodman%NULL%1]"
$ws.Cells.Item(8, 5).Value2 = "[Lauren P.%Jatt%NULL%1,             Alexander%Winnett%NULL%2,             Alexander%Winnett%NULL%0,             Christopher J.%Graber%NULL%2,             Christopher J.%Graber%NULL%0,             John%Vallone%NULL%1,             David O.%Beenhouwer%NULL%1,             Matthew Bidwell%Goetz%NULL%1]"
$ws.Cells.Item(9, 5).Value2 = "[Dan%Kelly%NULL%1,             Holly%Murphy%NULL%1,             Ravi%Vadlamudi%NULL%2,             Ravi%Vadlamudi%NULL%0,             Ruth%Kraut%NULL%1,             Kate%Dalessio%NULL%1,             Anurag N.%Malani%NULL%1,             Meghan%Glabach%NULL%1,             Juan Luis%Marquez%NULL%1]"
$ws.Cells.Item(10, 5).Value2 = "[Julia H.%Rogers%NULL%1,             Amy C.%Link%NULL%2,             Amy C.%Link%NULL%0,             Denise%McCulloch%NULL%1,             Elisabeth%Brandstetter%NULL%1,             Kira L.%Newman%NULL%2,             Kira L.%Newman%NULL%0,             Michael L.%Jackson%NULL%2,             Michael L.%Jackson%NULL%0,             James P.%Hughes%NULL%2,             James P.%Hughes%NULL%0,             Janet A.%Englund%NULL%2,             Janet A.%Englund%NULL%0,             Michael%Boeckh%NULL%2,             Michael%Boeckh%NULL%0,             Nancy%Sugg%NULL%2,             Misja%Ilcisin%NULL%1,             Thomas R.%Sibley%NULL%2,             Thomas R.%Sibley%NULL%0,             Kairsten%Fay%NULL%2,             Kairsten%Fay%NULL%0,             Jover%Lee%NULL%2,             Jover%Lee%NULL%0,             Peter%Han%NULL%2,             Peter%Han%NULL%0,             Melissa%Truong%NULL%1,             Matthew%Richardson%NULL%2,             Matthew%Richardson%NULL%0,             Deborah A.%Nickerson%NULL%2,             Deborah A.%Nickerson%NULL%0,             Lea M.%Starita%NULL%1,             Trevor%Bedford%NULL%2,             Trevor%Bedford%NULL%0,             Helen Y.%Chu%NULL%1,             NULL%NULL%NULL%0]"
$ws.Cells.Item(11, 5).Value2 = "[Matthew A.%Spinelli%xref no email%1,        Matthew D.%Hickey%xref no email%1,        David V.%Glidden%xref no email%0,        Janet Q.%Nguyen%xref no email%1,        Jon J.%Oskarsson%xref no email%1,        Diane%Havlir%xref no email%1,        Monica%Gandhi%xref no email%0]"
$ws.Cells.Item(12, 5).Value2 = "[Joan S.%Tucker%NULL%1,             Elizabeth J.%D'Amico%NULL%1,             Eric R.%Pedersen%NULL%1,             Rick%Garvey%NULL%1,             Anthony%Rodriguez%NULL%1,             David J.%Klein%NULL%1]"
$ws.Cells.Item(13, 5).Value2 = "[Hanna%Marquez%NULL%1,             Christian%Ramers%NULL%1,             Adam%Northrup%NULL%1,             Aaron%Tam%aaront@fhscd.org%1,             Jie%Liu%NULL%1,             Sarah%Rojas%NULL%1,             Stacey%Klaman%NULL%1,             Maureen%Khasira%NULL%1,             Jenan%Madbak%NULL%1,             Eva%Matthews%NULL%1,             Jeffrey%Norris%NULL%1,             Job%Godino%NULL%1]"
$ws.Cells.Item(14, 5).Value2 = "[Jane C%Yoon%NULL%1,             Martha P%Montgomery%NULL%2,             Martha P%Montgomery%NULL%0,             Ann M%Buff%NULL%1,             Andrew T%Boyd%NULL%1,             Calla%Jamison%NULL%0,             Alfonso%Hernandez%NULL%1,             Kristine%Schmit%NULL%0,             Sarita%Shah%NULL%1,             Sophia%Ajoku%NULL%1,             David P%Holland%NULL%1,             Juliana%Prieto%NULL%1,             Sasha%Smith%NULL%1,             Mark A%Swancutt%NULL%1,             Kim%Turner%NULL%1,             Tom%Andrews%NULL%0,             Kevin%Flowers%NULL%1,             Alyssa%Wells%NULL%1,             Cathryn%Marchman%NULL%1,             Emaline%Laney%NULL%1,             Danae%Bixler%NULL%1,             Sean%Cavanaugh%NULL%1,             Nicole%Flowers%NULL%1,             Nicholas%Gaffga%NULL%1,             Jean Y%Ko%NULL%1,             Heather N%Paulin%NULL%1,             Mark K%Weng%NULL%1,             Emily%Mosites%NULL%0,             Sapna Bamrah%Morris%feu3@cdc.gov%0]"
$ws.Cells.Item(15, 5).Value2 = "[Elizabeth%Imbert%Elizabeth.imbert@ucsf.edu%1,             Patrick M%Kinley%NULL%1,             Ashley%Scarborough%NULL%0,             Caroline%Cawley%NULL%1,             Madeline%Sankaran%NULL%0,             Sarah N%Cox%NULL%1,             Margot%Kushel%NULL%1,             Juliet%Stoltey%NULL%1,             Stephanie%Cohen%NULL%0,             Jonathan D%Fuchs%NULL%1,             NULL%NULL%NULL%0]"
$ws.Cells.Item(16, 5).Value2 = "[Farrell A.%Tobolowsky%NULL%1,             Elysia%Gonzales%NULL%0,             Julie L.%Self%NULL%0,             Carol Y.%Rao%NULL%0,             Ryan%Keating%NULL%0,             Grace E.%Marx%NULL%0,             Temet M.%McMichael%NULL%0,             Margaret D.%Lukoff%NULL%1,             Jeffrey S.%Duchin%NULL%1,             Karin%Huster%NULL%0,             Jody%Rauch%NULL%0,             Hedda%McLendon%NULL%1,             Matthew%Hanson%NULL%1,             Dave%Nichols%NULL%1,             Sargis%Pogosjans%NULL%1,             Meaghan%Fagalde%NULL%1,             Jennifer%Lenahan%NULL%1,             Emily%Maier%NULL%1,             Holly%Whitney%NULL%1,             Nancy%Sugg%NULL%0,             Helen%Chu%NULL%1,             Julia%Rogers%NULL%1,             Emily%Mosites%NULL%0,             Meagan%Kay%NULL%0]"
$ws.Cells.Item(17, 5).Value2 = "[Benjamin F.%Henwood%xref no email%1,        Brian%Redline%xref no email%1,        Jack%Lahey%xref no email%1]"
$ws.Cells.Item(18, 5).Value2 = "[Rebecca%Karb%NULL%1,             Elizabeth%Samuels%NULL%1,             Rahul%Vanjani%NULL%1,             Catherine%Trimbur%NULL%1,             Anthony%Napoli%NULL%1]"
$ws.Cells.Item(19, 5).Value2 = "[Travis P.%Baggett%xref no email%0,        Melanie W.%Racine%xref no email%1,        Elizabeth%Lewis%xref no email%1,        Denise%De Las Nueces%xref no email%1,        James J.%O\u2019Connell%xref no email%1,        Barry%Bock%xref no email%1,        Jessie M.%Gaeta%xref no email%0]"
$ws.Cells.Item(20, 5).Value2 = "[Jacqueline M%Ferguson%Jacqueline.Ferguson@va.gov%1,             Josephine%Jacobs%NULL%2,             Josephine%Jacobs%NULL%0,             Maria%Yefimova%NULL%1,             Liberty%Greene%NULL%1,             Leonie%Heyworth%NULL%1,             Donna M%Zulman%NULL%1]"
$ws.Cells.Item(21, 5).Value2 = "[Lucie%Richard%NULL%1,             Richard%Booth%NULL%1,             Jennifer%Rayner%NULL%1,             Kristin K.%Clemens%NULL%1,             Cheryl%Forchuk%NULL%1,             Salimah Z.%Shariff%NULL%1]"
$ws.Cells.Item(22, 5).Value2 = "[ Linwei%Wang%null%1,              Huiting%Ma%null%1,              Kristy C.Y.%Yiu%null%1,              Andrew%Calzavara%null%1,              David%Landsman%null%1,              Linh%Luong%null%1,              Adrienne K.%Chan%null%1,              Rafal%Kustra%null%1,              Jeffrey C.%Kwong%null%1,              Marie-Claude%Boily%null%1,              Stephen%Hwang%null%1,              Sharon%Straus%null%1,              Stefan D.%Baral%null%1,              Sharmistha%Mishra%null%1,             Linwei%Wang%null%1,             Huiting%Ma%null%1,             Kristy C.Y.%Yiu%null%1,             Andrew%Calzavara%null%1,             David%Landsman%null%1,             Linh%Luong%null%1,             Adrienne K.%Chan%null%1,             Rafal%Kustra%null%1,             Jeffrey C.%Kwong%null%1,             Marie-Claude%Boily%null%1,             Stephen%Hwang%null%1,             Sharon%Straus%null%1,             Stefan D.%Baral%null%1,             Sharmistha%Mishra%null%1]"
$ws.Cells.Item(23, 5).Value2 = "[Tran Duc Anh%Ly%NULL%0,             Van Thuan%Hoang%NULL%2,             Ndiaw%Goumballa%NULL%2,             Meriem%Louni%NULL%2,             Naomie%Canard%NULL%2,             Thi Loi%Dao%NULL%2,             Hacene%Medkour%NULL%2,             Audrey%Borg%NULL%2,             Kevin%Bardy%NULL%2,             Véra%Esteves-Vieira%NULL%2,             Véronique%Filosa%NULL%2,             Bernard%Davoust%NULL%2,             Oleg%Mediannikov%NULL%2,             Pierre-Edouard%Fournier%NULL%2,             Didier%Raoult%NULL%2,             Philippe%Gautret%philippe.gautret@club-internet.fr%2]"
$ws.Cells.Item(24, 5).Value2 = "[Thomas%Roederer%NULL%1,             Bastien%Mollo%NULL%1,             Charline%Vincent%NULL%1,             Birgit%Nikolay%NULL%1,             Augusto E%Llosa%NULL%1,             Robin%Nesbitt%NULL%1,             Jessica%Vanhomwegen%NULL%1,             Thierry%Rose%NULL%1,             Sophie%Goyard%NULL%1,             François%Anna%NULL%1,             Corinne%Torre%NULL%1,             Emilie%Fourrey%NULL%1,             Erica%Simons%NULL%1,             William%Hennequin%NULL%1,             Clair%Mills%NULL%1,             Francisco J%Luquero%NULL%1]"
$ws.Cells.Item(25, 5).Value2 = "[Tran Duc Anh%Ly%NULL%0,             Nhu Ngoc%Nguyen%NULL%1,             Van Thuan%Hoang%NULL%0,             Ndiaw%Goumballa%NULL%0,             Meriem%Louni%NULL%0,             Naomie%Canard%NULL%0,             Thi Loi%Dao%NULL%0,             Hacene%Medkour%NULL%0,             Audrey%Borg%NULL%0,             Kevin%Bardy%NULL%0,             Véra%Esteves-Vieira%NULL%0,             Véronique%Filosa%NULL%0,             Bernard%Davoust%NULL%0,             Oleg%Mediannikov%NULL%0,             Pierre-Edouard%Fournier%NULL%0,             Didier%Raoult%NULL%0,             Philippe%Gautret%NULL%0]"
$ws.Cells.Item(26, 5).Value2 = "[Massimo%Ralli%NULL%1,             Aldo%Morrone%NULL%1,             Andrea%Arcangeli%NULL%1,             Lucia%Ercoli%NULL%1]"
$ws.Cells.Item(27, 5).Value2 = "[Loïc%Schrooyen%NULL%1,             Marc%Delforge%NULL%1,             Faustine%Lebout%NULL%1,             Thibaut%Vanbaelen%NULL%1,             Amaryl%Lecompte%NULL%1,             Nicolas%Dauby%NULL%1]"
$ws.Cells.Item(28, 5).Value2 = "[Carmen%Martin%NULL%1,             Pilar%Andrés%NULL%2,             Pilar%Andrés%NULL%0,             Alberto%Bullón%NULL%1,             José Luis%Villegas%NULL%2,             José Luis%Villegas%NULL%0,             Javier Ignacio%de la Iglesia-Larrad%NULL%1,             Berta%Bote%NULL%1,             Nieves%Prieto%NULL%1,             Carlos%Roncero%NULL%1]"
